# BUG: Change "Chronic Renal Failure" to "Chronic Kidney Disease"
#
# The cause-of-death lookup table contains two cells referencing the old
# "Chronic Renal Failure" terminology (the cause-name cell and the
# ICD-10 description cell for code N18). Both need to read
# "Chronic Kidney Disease" instead.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdReplaceAll = 2, wdFindContinue = 1
$find.Execute(
    "Chronic Renal Failure",  # FindText
    $true,                    # MatchCase
    $true,                    # MatchWholeWord
    $false,                   # MatchWildcards
    $false,                   # MatchSoundsLike
    $false,                   # MatchAllWordForms
    $true,                    # Forward
    1,                        # Wrap -> wdFindContinue
    $false,                   # Format
    "Chronic Kidney Disease", # ReplaceWith
    2                         # Replace -> wdReplaceAll
)
